$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its values as exact text (avoids numeric
# auto-conversion/precision loss for values like "1.000" or "6.440").
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.191.90'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '1.902.06'
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '308.03'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = '0.5213'
$ws.Range('E7').Value = '  +1.21%  '
$ws.Range('D8').Value = '0.3766'
$ws.Range('E8').Value = '  +0.80%  '
$ws.Range('D9').Value = '0.07276'
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('D10').Value = '21.14'
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('D11').Value = '0.9050'
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').Value = '0.08289'
$ws.Range('E12').Value = '  +8.43%  '
$ws.Range('D13').Value = '1.907.89'
$ws.Range('E13').Value = '  +1.29%  '
$ws.Range('D14').Value = '96.71'
$ws.Range('E14').Value = '  +3.25%  '
$ws.Range('D15').Value = '5.291'
$ws.Range('E15').Value = '  +1.23%  '
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = '0.000008673'
$ws.Range('E17').Value = '  +2.38%  '
$ws.Range('D18').Value = '14.57'
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('D19').Value = '0.9998'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').Value = '27.228.81'
$ws.Range('E20').Value = '  +0.86%  '
$ws.Range('D21').Value = '5.093'
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('D22').Value = '2.157.67'
$ws.Range('E22').Value = '  +2.11%  '
$ws.Range('D23').Value = '10.65'
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').Value = '6.440'
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('D25').Value = '2.322'
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('D26').Value = '146.41'
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').Value = '1.747'
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '18.22'
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('D29').Value = '115.14'
$ws.Range('E29').Value = '  +1.14%  '
$ws.Range('D30').Value = '4.837'
$ws.Range('E30').Value = '  +1.40%  '
$ws.Range('D31').Value = '4.899'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').Value = '0.09274'
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('D33').Value = '0.05084'
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('D34').Value = '0.7989'
$ws.Range('E34').Value = '  +4.45%  '
$ws.Range('E35').Value = '  +1.46%  '
$ws.Range('D36').Value = '3.419'
$ws.Range('E36').Value = '  +4.90%  '
$ws.Range('D37').Value = '2.946'
$ws.Range('E37').Value = '  -1.39%  '
$ws.Range('D38').Value = '2.592'
$ws.Range('E38').Value = '  +0.46%  '
$ws.Range('D39').Value = '0.5726'
$ws.Range('E39').Value = '  +2.96%  '
$ws.Range('D40').Value = '0.02001'
$ws.Range('E40').Value = '  +0.94%  '
$ws.Range('D41').Value = '1.077'
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('D42').Value = '9.023'
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('D43').Value = '6.591'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').Value = '117.15'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('D45').Value = '0.1519'
$ws.Range('E45').Value = '  +1.46%  '
$ws.Range('D46').Value = '0.4867'
$ws.Range('E46').Value = '  +1.39%  '
$ws.Range('D47').Value = '1.000'
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('D48').Value = '10.14'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('E49').Value = '  +2.09%  '
$ws.Range('D50').Value = '37.72'
$ws.Range('D51').Value = '64.04'
$ws.Range('E51').Value = '  +0.50%  '
